# [rdbms] - remove double-commit scenario to avoid unnecessary error condition.
# - allow for non-standard or non-CRUD statements.
#
# The hidden "#system" sheet drives several named ranges (base, external, web, ...)
# that back the drop-down / autocomplete lists used elsewhere in the workbook.
# This change:
#   1. Removes the "clearVariables(variables)" entry from the "base" list
#      (column F), shifting everything below it up by one row.
#   2. Adds a "terminate(programName)" entry to the "external" list (column J).
#   3. Adds "saveSelectedText(var,locator)" and "saveSelectedValue(var,locator)"
#      entries to the "web" list (column Z), in their alphabetically-sorted spot.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- 1. column F ("base"): drop F19 ("clearVariables(variables)"), shift F20:F40 up ---
for ($r = 19; $r -le 39; $r++) {
    $ws.Cells.Item($r, 6).Value2 = $ws.Cells.Item($r + 1, 6).Value2
}
$ws.Cells.Item(40, 6).ClearContents()

# --- 2. column J ("external"): append "terminate(programName)" at J6 ---
$ws.Cells.Item(6, 10).Value2 = "terminate(programName)"

# --- 3. column Z ("web"): insert two new rows at Z99/Z100, shifting Z99:Z135 down to Z101:Z137 ---
for ($r = 135; $r -ge 99; $r--) {
    $ws.Cells.Item($r + 2, 26).Value2 = $ws.Cells.Item($r, 26).Value2
}
$ws.Cells.Item(99, 26).Value2 = "saveSelectedText(var,locator)"
$ws.Cells.Item(100, 26).Value2 = "saveSelectedValue(var,locator)"

# --- 4. keep the workbook-level defined names in sync with the resized ranges ---
$wb.Names.Item("base").RefersTo = "='#system'!`$F`$2:`$F`$39"
$wb.Names.Item("external").RefersTo = "='#system'!`$J`$2:`$J`$6"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$137"
